# Append a new案件 row at the top of the list, plus two more new rows
# further down the table, shifting the existing rows as needed, and
# refresh the "取得日時" timestamp on every data row to the new run time.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$newTimestamp = "2026-02-11 18:59:45"

# 1) Insert a brand-new top entry (becomes row 2); everything below
#    shifts down by one row.
$ws.Rows.Item(2).Insert()

# 2) Insert a new entry right before the "salesforce..." row (currently
#    row 10 after step 1), which becomes the new row 10.
$ws.Rows.Item(10).Insert()

# 3) Insert a new entry right before "【急募】ドメイン接続業務..." (currently
#    row 12 after step 2), which becomes the new row 12.
$ws.Rows.Item(12).Insert()

# This runtime does not re-target the worksheet's <hyperlinks> entries
# when rows are inserted (unlike real Excel), so drop every hyperlink
# now and rebuild the full set below once all rows are in their final
# positions.
$ws.Cells.Hyperlinks.Delete()

# Data for the three freshly-inserted rows.
$newRows = @(
    @{ Row = 2;  B = "【AI活用】市場調査・競合分析自動化機能の構築依頼"; C = "システム開発"; D = "50,000 円 ~ 100,000 円 / 固定"; E = "期限情報なし"; F = "https://www.lancers.jp/work/detail/5489981"; G = 388; H = "🔥AI,Ai ◆自動化" },
    @{ Row = 10; B = "初回 Bubbleを活用したマッチングシステム開発エンジニア募集"; C = "システム開発"; D = "5,000 円 ~ 10,000 円 / 固定"; E = "期限情報なし"; F = "https://www.lancers.jp/work/detail/5489949"; G = 110; H = "◆開発,システム開発" },
    @{ Row = 12; B = "プロジェクトマネジメント"; C = "システム開発"; D = "500,000 円 ~ 1,000,000 円 / 固定"; E = "期限情報なし"; F = "https://www.lancers.jp/work/detail/5490062"; G = 25; H = "" }
)

foreach ($item in $newRows) {
    $r = $item.Row
    $ws.Range("A$r").Value = $newTimestamp
    $ws.Range("B$r").Value = $item.B
    $ws.Range("C$r").Value = $item.C
    $ws.Range("D$r").Value = $item.D
    $ws.Range("E$r").Value = $item.E
    $ws.Range("F$r").Value = $item.F
    $ws.Range("G$r").Value = $item.G
    if ($item.H -ne "") {
        $ws.Range("H$r").Value = $item.H
    }
}

# 4) Refresh the "取得日時" timestamp for all of the rows that were not
#    freshly inserted above (their content is otherwise unchanged).
$unchangedRows = @(3, 4, 5, 6, 7, 8, 9, 11, 13, 14)
foreach ($r in $unchangedRows) {
    $ws.Range("A$r").Value = $newTimestamp
}

# 5) Rebuild every hyperlink (F2:F14) now that all rows are in their
#    final resting place, and make sure the Hyperlink cell style is
#    (re)applied.
$urls = @{
    2  = "https://www.lancers.jp/work/detail/5489981"
    3  = "https://www.lancers.jp/work/detail/5455098"
    4  = "https://www.lancers.jp/work/detail/5445159"
    5  = "https://www.lancers.jp/work/detail/5445154"
    6  = "https://www.lancers.jp/work/detail/5489818"
    7  = "https://www.lancers.jp/work/detail/5489711"
    8  = "https://www.lancers.jp/work/detail/5489911"
    9  = "https://www.lancers.jp/work/detail/5489608"
    10 = "https://www.lancers.jp/work/detail/5489949"
    11 = "https://www.lancers.jp/work/detail/5489898"
    12 = "https://www.lancers.jp/work/detail/5490062"
    13 = "https://www.lancers.jp/work/detail/5489674"
    14 = "https://www.lancers.jp/work/detail/5489636"
}

foreach ($r in @(2, 3, 4, 5, 6, 7, 8, 9, 10, 11, 12, 13, 14)) {
    $url = $urls[$r]
    $ws.Range("F$r").Value = $url
    $ws.Hyperlinks.Add($ws.Range("F$r"), $url)
    $ws.Range("F$r").Style = "Hyperlink"
}
